$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Rename the third column header to uppercase
$ws.Range("C1").Value = "LUS TSHAJ TAWM"

# Expand the table to include a new row
$tbl.Resize($ws.Range("A1:C3"))

# Copy the formatting of row 2 down into the new row 3
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)

# Fill in the new row's data
$ws.Range("A3").Value = [datetime]"2025-09-23"
$ws.Range("B3").Value = "Song practice this coming Thursday 9/25 will consist for both Sunday worship and 40 years anniversary."
$ws.Range("C3").Value = "Hnub Thursday, 9/25 no peb yuav kawm nkauj rau lub Sunday thiab rau lub 40 xyoo."

$ws.Rows.Item(3).RowHeight = 30

$ws.Range("C3").Select()
